$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'43.655.21"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  -1.00%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'2.302.92"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +2.17%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.10%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').Value = "'267.56"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  -1.79%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'93.94"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +7.17%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.44%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D9').Value = "'0.619"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.45%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'44.59"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -2.22%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.73%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'8.06"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +5.40%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('D13').Value = "'0.105"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.47%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').Value = "'2.644.01"
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').Value = "'15.37"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.45%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  +7.29%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'2.299.40"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +1.15%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'43.622.51"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  -0.90%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  +3.71%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('E20').Value = "'  +5.63%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').Value = "'  +0.88%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('E22').Value = "'  -4.60%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').Value = "'237.57"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +1.49%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'9.59"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +7.78%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'0.999"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.08%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').Value = "'11.31"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +3.94%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'2.49"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -0.82%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'3.40"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  -4.16%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'2.29"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.98%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'38.54"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  -3.98%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'22.28"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +6.51%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('E32').Value = "'  -1.95%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'0.0893"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.45%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'5.49"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +1.84%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').Value = "'  +1.62%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').Value = "'  -2.60%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').Value = "'0.0354"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.25%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'4.40"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +0.56%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  -2.03%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.233"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +14.42%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').Value = "'  +4.54%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').Value = "'  +19.80%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'12.09"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -4.75%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'5.43"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.36%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'61.78"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -4.39%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').Value = "'9.05"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +6.11%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').Value = "'  +3.39%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'100.48"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  -0.48%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').Value = "'  +0.41%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'2.523.21"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +1.98%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').Value = "'  -2.61%  "
$ws.Range('E51').Style = 'Normal'
